$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.985.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.04"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +17.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.07%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +16.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.580.76"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.825"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.86%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.241.73"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.994.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000106"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.83%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.35"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.77"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.96%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.88%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0912"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.35%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.14%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0351"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.07%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.66%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +27.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.253"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +26.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.55"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.54%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.46"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.68"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.02%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.11%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.436"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.460.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.82%  "
